$d = $word.ActiveDocument

# Colour constants (Word OLE_COLOR encoding):
#   wdColorAutomatic                      -> -16777216
#   Accent 3, Darker 50% (themeShade 80)  -> -704610049
#     (top byte 0xD6 = 0xD0 + themeColorIndex(accent3=6); next byte = themeShade
#      0x80; low byte 0xFF = "tint not set")
$wdColorAutomatic = -16777216
$accent3Dark50   = -704610049

# Paragraph 1: "What happens when you write this?" -> accent3 / darker 50%
$d.Paragraphs(1).Range.Font.Color = $accent3Dark50

# Paragraph 2: empty paragraph -> accent3 / darker 50%
$d.Paragraphs(2).Range.Font.Color = $accent3Dark50

# Paragraph 3: "Twinkle();" -> automatic
$d.Paragraphs(3).Range.Font.Color = $wdColorAutomatic

# Paragraph 4: empty paragraph -> automatic
$d.Paragraphs(4).Range.Font.Color = $wdColorAutomatic

# Paragraph 5: "As always, make sure you type it right, or the computer will
# not understand!" -> accent3 / darker 50%, and the run is split in two
# (around the existing _GoBack bookmark) right after "As always, make sure".
$p5 = $d.Paragraphs(5)
$full = $p5.Range
$paraStart = $full.Start
$paraTextEnd = $full.End - 1   # exclude the trailing paragraph mark
$splitAt = $paraStart + 20     # right after "As always, make sure"

$firstHalf = $d.Range($paraStart, $splitAt)
$secondHalf = $d.Range($splitAt, $paraTextEnd)

$firstHalf.Font.Color = $accent3Dark50
$secondHalf.Font.Color = $accent3Dark50
